# Update cryptocurrency price & volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.143.28"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "3.319.34"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.40"
$ws.Range("E5").Value = "  +2.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.89"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.312.66"
$ws.Range("E8").Value = "  +0.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.574"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.177"
$ws.Range("E10").Value = "  -4.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.574"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.50"
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000265"
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("D14").Value = "3.852.57"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.48"
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "615.12"
$ws.Range("E16").Value = "  -3.79%  "
$ws.Range("D17").Value = "66.170.36"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.87"
$ws.Range("E19").Value = "  -1.89%  "
$ws.Range("D20").Value = "3.321.09"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.98"
$ws.Range("E21").Value = "  -3.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.895"
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.20"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "99.85"
$ws.Range("E24").Value = "  -4.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.98"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.00"
$ws.Range("E26").Value = "  +1.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.72"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.50"
$ws.Range("E28").Value = "  -1.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "31.09"
$ws.Range("E29").Value = "  +2.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.49"
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.50"
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.76"
$ws.Range("E32").Value = "  -5.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "559.71"
$ws.Range("E33").Value = "  +4.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.89"
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("D35").Value = "3.832.03"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.05"
$ws.Range("E38").Value = "  -2.57%  "
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "32.77"
$ws.Range("E40").Value = "  -3.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.15"
$ws.Range("E41").Value = "  -4.08%  "
$ws.Range("D42").Value = "0.0₃0683"
$ws.Range("E42").Value = "  -7.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.60"
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("E44").Value = "  +4.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.334"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0408"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.07"
$ws.Range("E47").Value = "  -5.34%  "
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.52"
$ws.Range("E50").Value = "  -3.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.23"
$ws.Range("E51").Value = "  +5.30%  "
